$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.386.53"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.609.91"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0608"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.245"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.839.32"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "1.609.59"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.505"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "26.428.28"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.68%  "
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0494"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "1.451.26"
$ws.Range("E33").Value = "  +8.62%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.565"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.830"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "1.750.59"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.761"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0501"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.92%  "
